$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" '30.460.91'
Set-TextValue "E2" '  -0.67%  '
Set-TextValue "D3" '2.094.55'
Set-TextValue "E3" '  -0.64%  '
Set-TextValue "E4" '  -0.18%  '
Set-TextValue "D5" '329.36'
Set-TextValue "E6" '  -0.08%  '
Set-TextValue "D7" '0.5204'
Set-TextValue "E7" '  +0.11%  '
Set-TextValue "D8" '0.4340'
Set-TextValue "E8" '  -0.55%  '
Set-TextValue "D9" '52.68'
Set-TextValue "E9" '  +16.47%  '
Set-TextValue "D10" '0.08835'
Set-TextValue "E10" '  -1.40%  '
Set-TextValue "D11" '1.154'
Set-TextValue "E11" '  -1.37%  '
Set-TextValue "D12" '24.34'
Set-TextValue "E12" '  -1.85%  '
Set-TextValue "D13" '2.086.41'
Set-TextValue "E13" '  -1.89%  '
Set-TextValue "D14" '6.674'
Set-TextValue "E14" '  -1.42%  '
Set-TextValue "D15" '7.667'
Set-TextValue "E15" '  +0.06%  '
Set-TextValue "D16" '95.73'
Set-TextValue "E16" '  -1.75%  '
Set-TextValue "D17" '1.003'
Set-TextValue "E17" '  -0.02%  '
Set-TextValue "D18" '0.00001118'
Set-TextValue "E18" '  -1.62%  '
Set-TextValue "D19" '0.06580'
Set-TextValue "E19" '  -0.49%  '
Set-TextValue "D20" '19.19'
Set-TextValue "E20" '  +0.29%  '
Set-TextValue "E21" '  -0.06%  '
Set-TextValue "D22" '6.262'
Set-TextValue "E22" '  -2.37%  '
Set-TextValue "D23" '30.500.97'
Set-TextValue "E23" '  -1.19%  '
Set-TextValue "D24" '12.17'
Set-TextValue "E24" '  +1.76%  '
Set-TextValue "D25" '2.337'
Set-TextValue "E25" '  +3.37%  '
Set-TextValue "D26" '2.334.51'
Set-TextValue "E26" '  -1.74%  '
Set-TextValue "E27" '  -2.88%  '
Set-TextValue "D28" '2.584'
Set-TextValue "E28" '  +1.52%  '
Set-TextValue "D29" '162.41'
Set-TextValue "E29" '  -0.54%  '
Set-TextValue "D30" '131.50'
Set-TextValue "E30" '  -1.77%  '
Set-TextValue "D31" '1.191'
Set-TextValue "E31" '  +1.14%  '
Set-TextValue "D32" '0.1067'
Set-TextValue "E32" '  -0.11%  '
Set-TextValue "D33" '1.669'
Set-TextValue "E33" '  +10.19%  '
Set-TextValue "D34" '6.123'
Set-TextValue "E34" '  -1.09%  '
Set-TextValue "E35" '  -0.66%  '
Set-TextValue "D36" '10.09'
Set-TextValue "E36" '  +6.14%  '
Set-TextValue "D37" '0.02571'
Set-TextValue "E37" '  -0.23%  '
Set-TextValue "D38" '0.06791'
Set-TextValue "E38" '  +0.48%  '
Set-TextValue "D39" '12.70'
Set-TextValue "E39" '  -0.02%  '
Set-TextValue "D40" '5.445'
Set-TextValue "E40" '  -2.93%  '
Set-TextValue "D41" '0.2261'
Set-TextValue "E41" '  +0.80%  '
Set-TextValue "D42" '0.6902'
Set-TextValue "E42" '  +1.97%  '
Set-TextValue "E43" '  +1.27%  '
Set-TextValue "E44" '  -0.03%  '
Set-TextValue "D45" '0.6363'
Set-TextValue "E45" '  +1.39%  '
Set-TextValue "D46" '13.97'
Set-TextValue "E46" '  -2.67%  '
Set-TextValue "D47" '2.200'
Set-TextValue "E47" '  -1.89%  '
Set-TextValue "E48" '  -0.90%  '
Set-TextValue "D49" '1.236'
Set-TextValue "E49" '  +10.63%  '
Set-TextValue "D50" '1.239'
Set-TextValue "E50" '  -2.50%  '
Set-TextValue "D51" '81.80'
Set-TextValue "E51" '  -1.37%  '
